# String_ppt.pptx - "String Notes added for More Topics."
#
# Slide 1 has a comparison table contrasting creating a String with and
# without the `new` keyword. The two header cells had their labels
# swapped/corrected:
#   Cell(1,1): "new KeyWord()"            -> "Without new KeyWord()"
#   Cell(1,2): "Without new() keyword"    -> "new() keyword"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the table shape on the slide.
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $tableShape = $shp
    }
}

$tbl = $tableShape.Table

# Header-row, column 1: only the first run's text changes ("new " -> "Without new "),
# the remaining runs ("KeyWord", "()") stay untouched.
$cell1 = $tbl.Cell(1, 1)
$cell1.Shape.TextFrame.TextRange.Text = "Without new "

# Header-row, column 2: the single run's text changes entirely.
$cell2 = $tbl.Cell(1, 2)
$cell2.Shape.TextFrame.TextRange.Text = "new() keyword"
